# Auto-generated script applying odds updates per diff
$ws = $excel.ActiveWorkbook.ActiveSheet

# Row 12
$ws.Range("K12").Value = 4.4
$ws.Range("L12").Value = 1.65
$ws.Range("M12").Value = 2.12
$ws.Range("N12").Value = 2.87
$ws.Range("R12").Value = 2.15
$ws.Range("S12").Value = 1.62
$ws.Range("U12").Value = 18
$ws.Range("Z12").Value = 4.4
$ws.Range("AE12").Value = 5.6
$ws.Range("AF12").Value = 11
$ws.Range("AG12").Value = 10.25

# Row 13
$ws.Range("G13").Value = 2.2
$ws.Range("H13").Value = 2.62
$ws.Range("J13").Value = 1.15
$ws.Range("K13").Value = 4.75
$ws.Range("L13").Value = 1.6
$ws.Range("M13").Value = 2.22
$ws.Range("N13").Value = 2.72
$ws.Range("O13").Value = 1.4
$ws.Range("R13").Value = 2.15
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 5.3
$ws.Range("V13").Value = 9.5
$ws.Range("W13").Value = 23
$ws.Range("X13").Value = 23
$ws.Range("Y13").Value = 45
$ws.Range("Z13").Value = 4.75
$ws.Range("AB13").Value = 18
$ws.Range("AC13").Value = 120
$ws.Range("AE13").Value = 8

# Row 17
$ws.Range("G17").Value = 1.47
$ws.Range("H17").Value = 3.8
$ws.Range("I17").Value = 6.1
$ws.Range("N17").Value = 1.91
$ws.Range("O17").Value = 1.7
$ws.Range("T17").Value = 4.9
$ws.Range("U17").Value = 5.3
$ws.Range("V17").Value = 7.1
$ws.Range("W17").Value = 8
$ws.Range("X17").Value = 10.75
$ws.Range("Y17").Value = 25
$ws.Range("Z17").Value = 9
$ws.Range("AA17").Value = 6.6
$ws.Range("AB17").Value = 16.5
$ws.Range("AC17").Value = 80
$ws.Range("AF17").Value = 29
$ws.Range("AG17").Value = 16.5
$ws.Range("AH17").Value = 100
$ws.Range("AI17").Value = 55

# Row 18
$ws.Range("G18").Value = 1.72
$ws.Range("I18").Value = 4.15
$ws.Range("N18").Value = 1.91
$ws.Range("O18").Value = 1.7
$ws.Range("P18").Value = 1.36
$ws.Range("Q18").Value = 2.52
$ws.Range("T18").Value = 5.5
$ws.Range("U18").Value = 6.5
$ws.Range("W18").Value = 10.75
$ws.Range("Y18").Value = 23
$ws.Range("Z18").Value = 9.25
$ws.Range("AB18").Value = 14
$ws.Range("AC18").Value = 65
$ws.Range("AD18").Value = 500
$ws.Range("AF18").Value = 18
$ws.Range("AG18").Value = 11.75
$ws.Range("AH18").Value = 50
$ws.Range("AI18").Value = 32
$ws.Range("AJ18").Value = 37

# Row 26
$ws.Range("G26").Value = 1.3
$ws.Range("H26").Value = 4.55
$ws.Range("V26").Value = 7.4
$ws.Range("W26").Value = 6.6
$ws.Range("Z26").Value = 11.25
$ws.Range("AA26").Value = 8
$ws.Range("AB26").Value = 19.5
$ws.Range("AC26").Value = 90
$ws.Range("AE26").Value = 16

# Row 31
$ws.Range("H31").Value = 3.6
$ws.Range("L31").Value = 1.25
$ws.Range("M31").Value = 3.25
$ws.Range("N31").Value = 1.75
$ws.Range("O31").Value = 1.87
$ws.Range("T31").Value = 7.6
$ws.Range("U31").Value = 8.75
$ws.Range("X31").Value = 14
$ws.Range("Z31").Value = 11.25
$ws.Range("AA31").Value = 7
$ws.Range("AC31").Value = 65
